$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A1:C6")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A6"), 0, 2)
$ws.Sort.SetRange($range)
$ws.Sort.Header = 1
$ws.Sort.Apply()
